$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.06557403506259453
$ws.Range("C2").Value = 0.02556334423893891
$ws.Range("D2").Value = 0.05288602039217949
$ws.Range("E2").Value = 0.132447272676044
$ws.Range("F2").Value = 0.133973739915765
$ws.Range("G2").Value = 0.1332061333258237
$ws.Range("H2").Value = 0.02367381100663764

$ws.Range("B3").Value = 0.02015083550910759
$ws.Range("C3").Value = 0.006098061199845224
$ws.Range("D3").Value = 0.0008833148404090862
$ws.Range("E3").Value = 0.05604640154679212
$ws.Range("F3").Value = 0.1060689790631738
$ws.Range("G3").Value = 0.07334016759995782
$ws.Range("H3").Value = 0.009454616866257787
